$d = $word.ActiveDocument

$pairs = @(
    @("0.8346", "0.1017"),
    @("0.8488", "0.6807"),
    @("0.5996", "0.2698"),
    @("0.9509", "0.9107"),
    @("0.0440", "0.0738"),
    @("0.3254", "0.3438"),
    @("0.9416", "0.0115"),
    @("0.3389", "0.1208"),
    @("0.6204", "0.0519"),
    @("0.6374", "0.9421"),
    @("0.8737", "0.7737"),
    @("0.0465", "0.5515"),
    @("0.6601", "0.6117"),
    @("0.1252", "0.7731"),
    @("0.5935", "0.6567"),
    @("0.5273", "0.9344"),
    @("0.4950", "0.7547"),
    @("0.3916", "0.4475"),
    @("0.4026", "0.0711"),
    @("0.9885", "0.4099"),
    @("0.9312", "0.4539"),
    @("0.7953", "0.9248"),
    @("0.4480", "0.7428"),
    @("0.8693", "0.3210"),
    @("0.7849", "0.0240"),
    @("0.8217", "0.7982"),
    @("0.6946", "0.6121"),
    @("0.7747", "0.4240"),
    @("0.6560", "0.4531"),
    @("0.7062", "0.7448"),
    @("0.8460", "0.2946"),
    @("0.0736", "0.2014"),
    @("0.1865", "0.4369"),
    @("0.1341", "0.6924"),
    @("0.9357", "0.3289"),
    @("0.2280", "0.6965"),
    @("0.9530", "0.9393"),
    @("0.8554", "0.7680"),
    @("0.9254", "0.1906"),
    @("0.4335", "0.8381"),
    @("0.6125", "0.7801"),
    @("0.5116", "0.9469"),
    @("0.8247", "0.8054"),
    @("0.8642", "0.6672"),
    @("0.1565", "0.2239"),
    @("0.7079", "0.7327"),
    @("0.6095", "0.3814"),
    @("0.4267", "0.0062"),
    @("0.8173", "0.0675"),
    @("0.1736", "0.4721"),
    @("0.8754", "0.6788"),
    @("0.7660", "0.2520"),
    @("0.9009", "0.5230"),
    @("0.1527", "0.6811"),
    @("0.4818", "0.0610"),
    @("0.2637", "0.1506")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
